$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM values for rows 2-7 (columns E-T), per the new script output.
$data = @{
    2 = @{ E=3; F=1; G=0.1292686666666667; H=0.387806; M=1.090710333333333; N=3.272131; O=0.0488470045579656; P=0.0488470045579656; Q=0.1409946705095556; R=1.268952034586;  S=0.0488470045579656; T=0.0488470045579656 }
    3 = @{ E=3; F=1; G=0.1292686666666667; H=0.387806;                                                          O=0.7616320856558244; P=0.7616320856558244; Q=2.198416585383778;  R=19.785749268454;    S=0.7616320856558244; T=0.7616320856558244 }
    4 = @{ E=3; F=1; G=0.1292686666666667; H=0.387806; M=0.740281;           N=2.220843;  O=0.03315317392351528; P=0.03315317392351528; Q=0.09569513782866666; R=0.861256240458;    S=0.03315317392351528; T=0.03315317392351528 }
    5 = @{ E=3; F=1; G=0.1292686666666667; H=0.387806; M=2.784013333333333;  N=8.352039999999999; O=0.1246808688124989; P=0.1246808688124989; Q=0.3598856915822222; R=3.238971224239999; S=0.1246808688124989; T=0.1246808688124989 }
    6 = @{ E=3; F=1; G=0.1292686666666667; H=0.387806; M=0.2710316666666667; N=0.813095;  O=0.01213803945228936; P=0.01213803945228936; Q=0.03503590217444445; R=0.31532311957;     S=0.01213803945228936; T=0.01213803945228936 }
    7 = @{ E=3; F=1; G=0.1292686666666667; H=0.387806; M=0.436508;           N=1.309524;  O=0.01954882759790648; P=0.01954882759790648; Q=0.05642680714933333; R=0.5078412643439999; S=0.01954882759790648; T=0.01954882759790648 }
}

foreach ($rowNum in $data.Keys) {
    $rowValues = $data[$rowNum]
    foreach ($col in $rowValues.Keys) {
        $addr = "$col$rowNum"
        $ws.Range($addr).Value = $rowValues[$col]
    }
}
